$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 338.2
$ws.Range("I4").Value = 385.25
$ws.Range("J4").Value = 150
$ws.Range("K4").Value = 385.25
$ws.Range("L4").Value = 150
$ws.Range("M4").Value = -271.25
$ws.Range("N4").Value = -378

# Row 9
$ws.Range("H9").Value = 7936895
$ws.Range("I9").Value = 657
$ws.Range("K9").Value = 657
$ws.Range("M9").Value = -488

# Row 107
$ws.Range("H107").Value = 9657.546
$ws.Range("I107").Value = 11548.556
$ws.Range("K107").Value = 11548.556
$ws.Range("M107").Value = -9628.556

# Row 117
$ws.Range("H117").Value = 371022500
$ws.Range("J117").Value = 371022500
$ws.Range("L117").Value = 371022500
$ws.Range("N117").Value = -371031678

$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 910.5
$ws.Range("I5").Value = 275
$ws.Range("J5").Value = 1546
$ws.Range("K5").Value = 275
$ws.Range("L5").Value = 1546
$ws.Range("M5").Value = -163
$ws.Range("N5").Value = -1770

# Row 32
$ws.Range("H32").Value = 2189.9265
$ws.Range("I32").Value = 2160.2307
$ws.Range("K32").Value = 2160.2307
$ws.Range("M32").Value = -1873.2307

# Row 43
$ws.Range("H43").Value = 17453.1
$ws.Range("I43").Value = 16578.666
$ws.Range("J43").Value = 17827.857
$ws.Range("K43").Value = 16578.666
$ws.Range("L43").Value = 17827.857
$ws.Range("M43").Value = -16265.666
$ws.Range("N43").Value = -18453.857

# Row 45
$ws.Range("H45").Value = 87992.8
$ws.Range("I45").Value = 120676.555
$ws.Range("K45").Value = 120676.555
$ws.Range("M45").Value = -120299.555

# Row 74
$ws.Range("H74").Value = 10669900
$ws.Range("I74").Value = 13196566
$ws.Range("K74").Value = 13196566
$ws.Range("M74").Value = -13195692

# Row 77
$ws.Range("H77").Value = 10669900
$ws.Range("I77").Value = 13196566
$ws.Range("K77").Value = 65982830
$ws.Range("M77").Value = -65978462

# Row 122
$ws.Range("H122").Value = 382142.53
$ws.Range("I122").Value = 2862.074
$ws.Range("K122").Value = 8586.222
$ws.Range("M122").Value = -6136.222

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 910.5
$ws.Range("I4").Value = 275
$ws.Range("J4").Value = 1546
$ws.Range("K4").Value = 275
$ws.Range("L4").Value = 1546
$ws.Range("M4").Value = -160
$ws.Range("N4").Value = -1776

# Row 80
$ws.Range("H80").Value = 586.4737
$ws.Range("I80").Value = 379.83334
$ws.Range("K80").Value = 379.83334
$ws.Range("M80").Value = 618.16666

# Row 83
$ws.Range("H83").Value = 586.4737
$ws.Range("I83").Value = 379.83334
$ws.Range("K83").Value = 1899.1667
$ws.Range("M83").Value = 3092.8333

# Row 88
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

# Row 91
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

# Row 105
$ws.Range("H105").Value = 45259.96
$ws.Range("I105").Value = 52216.453
$ws.Range("K105").Value = 52216.453
$ws.Range("M105").Value = -50469.453

# Row 134
$ws.Range("H134").Value = 5815
$ws.Range("I134").Value = 6012.6562
$ws.Range("J134").Value = 4233.75
$ws.Range("K134").Value = 18037.9686
$ws.Range("L134").Value = 12701.25
$ws.Range("M134").Value = -15502.9686
$ws.Range("N134").Value = -17771.25

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 10989685
$ws.Range("I22").Value = 780
$ws.Range("J22").Value = 15385246
$ws.Range("K22").Value = 780
$ws.Range("L22").Value = 15385246
$ws.Range("M22").Value = -430
$ws.Range("N22").Value = -15385946

# Row 31
$ws.Range("H31").Value = 2285.65
$ws.Range("J31").Value = 2559.2222
$ws.Range("L31").Value = 2559.2222
$ws.Range("N31").Value = -3149.2222

# Row 34
$ws.Range("H34").Value = 2285.65
$ws.Range("J34").Value = 2559.2222
$ws.Range("L34").Value = 2559.2222
$ws.Range("N34").Value = -2963.2222

# Row 50
$ws.Range("H50").Value = 29999.666
$ws.Range("J50").Value = 29999.666
$ws.Range("L50").Value = 29999.666
$ws.Range("N50").Value = -31249.666

# Row 58
$ws.Range("H58").Value = 6364.436
$ws.Range("I58").Value = 8719.632
$ws.Range("K58").Value = 8719.632
$ws.Range("M58").Value = -8516.632

# Row 60
$ws.Range("H60").Value = 23236
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 23236
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 23236
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -24258

# Row 122
$ws.Range("H122").Value = 2609.5
$ws.Range("I122").Value = 2146.3333
$ws.Range("J122").Value = 3999
$ws.Range("K122").Value = 6438.999899999999
$ws.Range("L122").Value = 11997
$ws.Range("M122").Value = -3988.999899999999
$ws.Range("N122").Value = -16897

# Row 132
$ws.Range("H132").Value = 19639992
$ws.Range("I132").Value = 25651528
$ws.Range("J132").Value = 102500
$ws.Range("K132").Value = 76954584
$ws.Range("L132").Value = 307500
$ws.Range("M132").Value = -76952054
$ws.Range("N132").Value = -312560

# Row 134
$ws.Range("H134").Value = 1650236
$ws.Range("I134").Value = 2724313.5
$ws.Range("J134").Value = 3317.1333
$ws.Range("K134").Value = 8172940.5
$ws.Range("L134").Value = 9951.3999
$ws.Range("M134").Value = -8170405.5
$ws.Range("N134").Value = -15021.3999

# Row 136
$ws.Range("H136").Value = 6364.436
$ws.Range("I136").Value = 8719.632
$ws.Range("K136").Value = 26158.896
$ws.Range("M136").Value = -23608.896

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 41374170
$ws.Range("I4").Value = 59755544
$ws.Range("K4").Value = 179266632
$ws.Range("M4").Value = -179266520

# Row 12
$ws.Range("H12").Value = 13.25
$ws.Range("I12").Value = 17
$ws.Range("K12").Value = 51
$ws.Range("M12").Value = 122

# Row 41
$ws.Range("H41").Value = 2057.45
$ws.Range("J41").Value = 353.07693
$ws.Range("L41").Value = 1059.23079
$ws.Range("N41").Value = -1735.23079

# Row 123
$ws.Range("H123").Value = 5419.8
$ws.Range("I123").Value = 4999.5
$ws.Range("J123").Value = 5700
$ws.Range("K123").Value = 14998.5
$ws.Range("L123").Value = 17100
$ws.Range("M123").Value = -12548.5
$ws.Range("N123").Value = -22000

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 4159.2163
$ws.Range("I132").Value = 3311.3428
$ws.Range("J132").Value = 18997
$ws.Range("K132").Value = 9934.028399999999
$ws.Range("L132").Value = 56991
$ws.Range("M132").Value = -7404.028399999999
$ws.Range("N132").Value = -62051

$ws = $wb.Worksheets.Item("LTW")
# Row 20
$ws.Range("H20").Value = 3562.4375
$ws.Range("I20").Value = 2300.3
$ws.Range("K20").Value = 2300.3
$ws.Range("M20").Value = -2074.3

# Row 68
$ws.Range("H68").Value = 4112.0625
$ws.Range("I68").Value = 1942.1428
$ws.Range("J68").Value = 5799.778
$ws.Range("K68").Value = 1942.1428
$ws.Range("L68").Value = 5799.778
$ws.Range("M68").Value = -1193.1428
$ws.Range("N68").Value = -7297.778

# Row 71
$ws.Range("H71").Value = 4112.0625
$ws.Range("I71").Value = 1942.1428
$ws.Range("J71").Value = 5799.778
$ws.Range("K71").Value = 9710.714
$ws.Range("L71").Value = 28998.89
$ws.Range("M71").Value = -5966.714
$ws.Range("N71").Value = -36486.89

# Row 132
$ws.Range("H132").Value = 12533.963
$ws.Range("I132").Value = 17607.53
$ws.Range("J132").Value = 3908.9
$ws.Range("K132").Value = 52822.59
$ws.Range("L132").Value = 11726.7
$ws.Range("M132").Value = -50292.59
$ws.Range("N132").Value = -16786.7

# Row 136
$ws.Range("H136").Value = 5923.9697
$ws.Range("I136").Value = 2456.4666
$ws.Range("J136").Value = 8813.556
$ws.Range("K136").Value = 7369.399800000001
$ws.Range("L136").Value = 26440.668
$ws.Range("M136").Value = -4819.399800000001
$ws.Range("N136").Value = -31540.668

$ws = $wb.Worksheets.Item("WVR")
# Row 41
$ws.Range("H41").Value = 13449.5
$ws.Range("J41").Value = 12539.4
$ws.Range("L41").Value = 12539.4
$ws.Range("N41").Value = -13319.4

# Row 62
$ws.Range("H62").Value = 196942.64
$ws.Range("I62").Value = 541800.4
$ws.Range("J62").Value = 5355
$ws.Range("K62").Value = 541800.4
$ws.Range("L62").Value = 5355
$ws.Range("M62").Value = -541176.4
$ws.Range("N62").Value = -6603

# Row 65
$ws.Range("H65").Value = 196942.64
$ws.Range("I65").Value = 541800.4
$ws.Range("J65").Value = 5355
$ws.Range("K65").Value = 2709002
$ws.Range("L65").Value = 26775
$ws.Range("M65").Value = -2705882
$ws.Range("N65").Value = -33015
